# Applies the "changed MP time limit and corrected error in fixed recourse
# data" commit to the random_recourse summary workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrected / updated raw data -----------------------------------------
# Row 9 (M=20,N=20,T=10,alpha=0.25): multi-cut PCP time + G column corrected
$ws.Range("F9").Value = "3018.2(9)"
$ws.Range("G9").Value = "1432.4(5)"
$ws.Range("H9").Value = 12
$ws.Range("J9").Value = 11

# Row 18 (M=50,N=100,T=50,alpha=0.25): MP time limit sentinel N/A -> "-", G corrected
$ws.Range("E18").Value = "- (10)"
$ws.Range("F18").Value = "- (10)"
$ws.Range("G18").Value = 1594.2
$ws.Range("H18").Value = "-"
$ws.Range("J18").Value = "-"

# Row 19 (M=50,N=100,T=50,alpha=0.5): MP time limit sentinel N/A -> "-", G corrected
$ws.Range("E19").Value = "- (10)"
$ws.Range("F19").Value = "- (10)"
$ws.Range("G19").Value = 1493
$ws.Range("H19").Value = "-"
$ws.Range("J19").Value = "-"

# Row 20 (M=50,N=100,T=50,alpha=0.75): MP time limit sentinel N/A -> "-", G corrected
$ws.Range("E20").Value = "- (10)"
$ws.Range("F20").Value = "2567.8(8)"
$ws.Range("G20").Value = 1189.3

# --- New column M: relative optimality-gap formulas ------------------------
# Standard rows use (E-G)/E against the row's own reformulation time.
$stdRows = 3,4,5,7,8,10,11,12,13,14,16,17
foreach ($r in $stdRows) {
    $ws.Range("M$r").Formula = "=(E$r-G$r)/E$r"
}

# Rows whose E-column baseline is text (N/A-style entries) use the
# corresponding numeric baseline instead of the E-column cell.
$ws.Range("M6").Formula = "=(740.6-G6)/740.6"
$ws.Range("M9").Formula = "=(140-1432)/140"
$ws.Range("M15").Formula = "=(1015-G15)/1015"
$ws.Range("M18").Formula = "=(3600-G18)/3600"
$ws.Range("M19").Formula = "=(3600-G19)/3600"
$ws.Range("M20").Formula = "=(3600-G20)/3600"

# Rows 18-20 previously had no formatting on column M (the row stopped at K);
# give them the same "0.00" numeric style already used by M3:M17.
$ws.Range("M18").NumberFormat = "0.00"
$ws.Range("M19").NumberFormat = "0.00"
$ws.Range("M20").NumberFormat = "0.00"

# --- Cursor position --------------------------------------------------------
$ws.Range("H37").Select() | Out-Null
